$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the "Total Expenses" row (currently row 31), shifting
# the totals/footer rows down by one.
$ws.Rows.Item(31).Insert()

# Fill in the new expense line: "Amazon Order 15" for Heat Shrink Tubing.
$ws.Range("A31").Value2 = "Amazon Order 15"
$ws.Range("B31").Value2 = 43191
$ws.Range("C31").Value2 = "Brian"
$ws.Range("D31").Value2 = "Amazon Order 15.pdf"
$ws.Range("E31").Value2 = 16.98
$ws.Range("F31").Value2 = "Heat Shrink Tubing"

# Add the hyperlink to the receipt PDF, then restore the cell to the shared
# "Hyperlink" cell style used by every other receipt cell (Hyperlinks.Add
# otherwise stamps direct font formatting onto the cell).
$ws.Hyperlinks.Add($ws.Range("D31"), "Amazon Order 15.pdf", "", "", "Amazon Order 15.pdf")
$ws.Range("D31").Style = "Hyperlink"

# Extend the Total Expenses sum to include the new row.
$ws.Range("E32").Formula = "=SUM(E2:E31)"

# Match the saved selection from the edit.
$ws.Range("F31").Select()
